$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 892056.9399999999
$ws.Range("I11").Value = 892056.9399999999
$ws.Range("K11").Value = 892056.9399999999
$ws.Range("M11").Value = -891916.9399999999
$ws.Range("H15").Value = 88.2
$ws.Range("I15").Value = 88.2
$ws.Range("K15").Value = 264.6
$ws.Range("M15").Value = -95.60000000000002
$ws.Range("H29").Value = 1255.4
$ws.Range("I29").Value = 569.25
$ws.Range("J29").Value = 4000
$ws.Range("K29").Value = 1707.75
$ws.Range("L29").Value = 12000
$ws.Range("M29").Value = -1426.75
$ws.Range("N29").Value = -12562
$ws.Range("H40").Value = 3054
$ws.Range("I40").Value = 1166.6666
$ws.Range("K40").Value = 1166.6666
$ws.Range("M40").Value = -991.6666
$ws.Range("H92").Value = 1256.0834
$ws.Range("I92").Value = 587
$ws.Range("J92").Value = 3263.3333
$ws.Range("K92").Value = 587
$ws.Range("L92").Value = 3263.3333
$ws.Range("M92").Value = 661
$ws.Range("N92").Value = -5759.3333
$ws.Range("H125").Value = 2648
$ws.Range("I125").Value = 3490
$ws.Range("J125").Value = 2086.6667
$ws.Range("K125").Value = 31410
$ws.Range("L125").Value = 18780.0003
$ws.Range("M125").Value = -28950
$ws.Range("N125").Value = -23700.0003
$ws.Range("H127").Value = 1770.9333
$ws.Range("I127").Value = 716.8
$ws.Range("J127").Value = 2298
$ws.Range("K127").Value = 2150.4
$ws.Range("L127").Value = 6894
$ws.Range("M127").Value = 2809.6
$ws.Range("N127").Value = -16814
$ws.Range("H137").Value = 2050.5293
$ws.Range("I137").Value = 1366.907
$ws.Range("J137").Value = 5725
$ws.Range("K137").Value = 4100.721
$ws.Range("L137").Value = 17175
$ws.Range("M137").Value = -1550.721
$ws.Range("N137").Value = -22275

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 7500
$ws.Range("I25").Value = 4000
$ws.Range("J25").Value = 11000
$ws.Range("K25").Value = 4000
$ws.Range("L25").Value = 11000
$ws.Range("M25").Value = -3598
$ws.Range("N25").Value = -11804
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("M30").ClearContents()
$ws.Range("N30").ClearContents()
$ws.Range("H61").Value = 1683.16
$ws.Range("I61").Value = 1037.4375
$ws.Range("J61").Value = 2831.111
$ws.Range("K61").Value = 1037.4375
$ws.Range("L61").Value = 2831.111
$ws.Range("M61").Value = -825.4375
$ws.Range("N61").Value = -3255.111
$ws.Range("H132").Value = 2377.9487
$ws.Range("I132").Value = 1167.8077
$ws.Range("K132").Value = 3503.4231
$ws.Range("M132").Value = -973.4231
$ws.Range("H136").Value = 1683.16
$ws.Range("I136").Value = 1037.4375
$ws.Range("J136").Value = 2831.111
$ws.Range("K136").Value = 3112.3125
$ws.Range("L136").Value = 8493.332999999999
$ws.Range("M136").Value = -562.3125
$ws.Range("N136").Value = -13593.333

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2871.3635
$ws.Range("I31").Value = 1317.4517
$ws.Range("J31").Value = 6576.846
$ws.Range("K31").Value = 1317.4517
$ws.Range("L31").Value = 6576.846
$ws.Range("M31").Value = -1022.4517
$ws.Range("N31").Value = -7166.846
$ws.Range("H34").Value = 2871.3635
$ws.Range("I34").Value = 1317.4517
$ws.Range("J34").Value = 6576.846
$ws.Range("K34").Value = 1317.4517
$ws.Range("L34").Value = 6576.846
$ws.Range("M34").Value = -1115.4517
$ws.Range("N34").Value = -6980.846
$ws.Range("H99").Value = 3153.318
$ws.Range("I99").Value = 1809.4117
$ws.Range("J99").Value = 7722.6
$ws.Range("K99").Value = 1809.4117
$ws.Range("L99").Value = 7722.6
$ws.Range("M99").Value = -311.4117000000001
$ws.Range("N99").Value = -10718.6
$ws.Range("H126").Value = 3153.318
$ws.Range("I126").Value = 1809.4117
$ws.Range("J126").Value = 7722.6
$ws.Range("K126").Value = 5428.2351
$ws.Range("L126").Value = 23167.8
$ws.Range("M126").Value = -2958.2351
$ws.Range("N126").Value = -28107.8
$ws.Range("H132").Value = 3562.7097
$ws.Range("I132").Value = 3554.5293
$ws.Range("J132").Value = 3572.6428
$ws.Range("K132").Value = 10663.5879
$ws.Range("L132").Value = 10717.9284
$ws.Range("M132").Value = -8133.5879
$ws.Range("N132").Value = -15777.9284

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 547.9
$ws.Range("I113").Value = 566.1177
$ws.Range("J113").Value = 524.0769
$ws.Range("K113").Value = 1698.3531
$ws.Range("L113").Value = 1572.2307
$ws.Range("M113").Value = 471.6469
$ws.Range("N113").Value = -5912.2307

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3956.1382
$ws.Range("I126").Value = 2884.2942
$ws.Range("J126").Value = 5227.3955
$ws.Range("K126").Value = 8652.882599999999
$ws.Range("L126").Value = 15682.1865
$ws.Range("M126").Value = -6182.882599999999
$ws.Range("N126").Value = -20622.1865
$ws.Range("H139").Value = 69980
$ws.Range("J139").Value = 69980
$ws.Range("L139").Value = 69980
$ws.Range("N139").Value = -80260

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7762.5
$ws.Range("I7").Value = 5266.6665
$ws.Range("J7").Value = 9260
$ws.Range("K7").Value = 5266.6665
$ws.Range("L7").Value = 9260
$ws.Range("M7").Value = -5154.6665
$ws.Range("N7").Value = -9484
$ws.Range("H68").Value = 889.32556
$ws.Range("I68").Value = 724.64105
$ws.Range("J68").Value = 2495
$ws.Range("K68").Value = 724.64105
$ws.Range("L68").Value = 2495
$ws.Range("M68").Value = 24.35895000000005
$ws.Range("N68").Value = -3993
$ws.Range("H71").Value = 889.32556
$ws.Range("I71").Value = 724.64105
$ws.Range("J71").Value = 2495
$ws.Range("K71").Value = 3623.20525
$ws.Range("L71").Value = 12475
$ws.Range("M71").Value = 120.79475
$ws.Range("N71").Value = -19963
$ws.Range("H126").Value = 7762.5
$ws.Range("I126").Value = 5266.6665
$ws.Range("J126").Value = 9260
$ws.Range("K126").Value = 15799.9995
$ws.Range("L126").Value = 27780
$ws.Range("M126").Value = -13329.9995
$ws.Range("N126").Value = -32720

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 9528917
$ws.Range("I132").Value = 10889.6
$ws.Range("K132").Value = 32668.8
$ws.Range("M132").Value = -30138.8
$ws.Range("H138").Value = 51422.727
$ws.Range("J138").Value = 51422.727
$ws.Range("L138").Value = 51422.727
$ws.Range("N138").Value = -61702.727
